# Insert a new data-collection column ("04dec2025") immediately to the left
# of the existing "26nov2025" column (column E) on both worksheets, shifting
# the old E:I columns to F:J, then populate the new column's header and
# per-empadronador counts.

$wb = $excel.ActiveWorkbook

# New header text for the inserted column.
$newHeader = "04dec2025"

# New values for the inserted column, keyed by row number (row 1 is the
# header row; rows 2-11 are the empadronador data rows).
$newValues = @{
    2  = 16
    3  = 15
    4  = 13
    5  = 12
    6  = 14
    7  = 16
    8  = 8
    9  = 12
    10 = 13
    11 = 17
}

foreach ($ws in $wb.Worksheets) {
    # Shift existing columns E:I to F:J, creating a blank column E.
    $ws.Columns("E:E").Insert()

    # Header cell (inherits the bold/border style from the insert, matching
    # the other header cells).
    $ws.Range("E1").Value = $newHeader

    # "crosstab" stores the counts as numbers; "annot" stores the very same
    # counts as text (matching how every other data column on that sheet is
    # typed). A leading apostrophe forces Excel to keep a numeric-looking
    # entry as text instead of silently coercing it back to a number.
    $isText = ($ws.Name -eq "annot")

    foreach ($row in 2..11) {
        $val = $newValues[$row]
        if ($isText) {
            $ws.Cells.Item($row, 5).Value = "'" + [string]$val
        } else {
            $ws.Cells.Item($row, 5).Value = $val
        }
    }
}
